$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two download-range cells whose values changed.
$ws.Range("G5").Value = "43:48"
$ws.Range("E16").Value = "43:48"
$ws.Range("F16").Value = "49:54"
$ws.Range("G16").Value = "55:60"

# Update the active-cell selection to match the new state.
$ws.Range("G6").Select()
